$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.007.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.502.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.94%  "
$ws.Range("E7").Value = "  +0.61%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.568"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.508.48"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("E10").Value = "  +3.25%  "
$ws.Range("E11").Value = "  -2.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.330"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.950.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.815.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.36"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.56%  "
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.506.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "322.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.69%  "
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.408"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.99%  "
$ws.Range("E26").Value = "  +1.01%  "
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "175.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0757"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.61%  "
$ws.Range("E31").Value = "  +1.84%  "
$ws.Range("E32").Value = "  +1.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.69%  "
$ws.Range("E35").Value = "  +0.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.22"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.92"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.32%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.61"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.44%  "
$ws.Range("B41").Value = "SuiNetwork"
$ws.Range("C41").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.817"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "275.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "131.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.591"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0940"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0508"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.07%  "
$ws.Range("E49").Value = "  +3.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.748.38"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.56%  "
